$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("buffer")

$data = @(
    @(5, "ANTIOQUIA", 1.5),
    @(8, "ATLÁNTICO", 1.5),
    @(11, "BOGOTÁ, D.C.", 1.5),
    @(13, "BOLÍVAR", 1.5),
    @(15, "BOYACÁ", 1.5),
    @(17, "CALDAS", 1.5),
    @(18, "CAQUETÁ", 1.5),
    @(19, "CAUCA", 1.5),
    @(20, "CESAR", 1.5),
    @(23, "CÓRDOBA", 1.5),
    @(25, "CUNDINAMARCA", 1.5),
    @(27, "CHOCÓ", 1.5),
    @(41, "HUILA", 1.5),
    @(44, "LA GUAJIRA", 1.5),
    @(47, "MAGDALENA", 1.5),
    @(50, "META", 1.5),
    @(52, "NARIÑO", 1.5),
    @(54, "NORTE DE SANTANDER", 1.5),
    @(63, "QUINDIO", 1.5),
    @(66, "RISARALDA", 1.5),
    @(68, "SANTANDER", 1.5),
    @(70, "SUCRE", 1.5),
    @(73, "TOLIMA", 1.5),
    @(76, "VALLE DEL CAUCA", 1.5),
    @(81, "ARAUCA", 1.5),
    @(85, "CASANARE", 1.5),
    @(86, "PUTUMAYO", 1.5),
    @(88, "ARCHIPIÉLAGO DE SAN ANDRÉS, PROVIDENCIA Y SANTA CATALINA", 1.5),
    @(91, "AMAZONAS", 1.5),
    @(94, "GUAINÍA", 1.5),
    @(95, "GUAVIARE", 1.5),
    @(97, "VAUPÉS", 1.5),
    @(99, "VICHADA", 1.5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
}

$ws.Range("C6").Select() | Out-Null